# Updated cryptos list on Fri Nov  3 14:09:01 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text that often looks numeric ("39.86", "0.615", ...).
# Assigning such a string straight to .Value lets Excel auto-convert it to a real number,
# which would silently drop meaningful trailing zeros (e.g. "4.60" -> 4.6). Flip the cell
# to Text format before writing so the value round-trips as a string exactly as typed,
# then restore the default "Normal" style so no stray formatting is left behind.

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.812.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.21%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6: XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.11%  "

# Row 8: Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.78%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.323"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.28%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0682"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.64%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -2.06%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.079.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.22%  "

# Row 13: Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "

# Rows 14 & 15: Polygon and WrappedEther swap ranking positions
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.836.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.18%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.668"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "

# Row 16: Polkadot
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.97%  "

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.788.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.29%  "

# Row 18: Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.51%  "

# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0782"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "

# Row 21: Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.69%  "

# Row 23: Dai
$ws.Range("E23").Value = "  -0.11%  "

# Row 24: Toncoin
$ws.Range("E24").Value = "  -0.30%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.53%  "

# Row 27: Stellar
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.123"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.29%  "

# Row 28: EthereumClassic
$ws.Range("E28").Value = "  -2.14%  "

# Row 29: PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.37%  "

# Row 30: BinanceUSD
$ws.Range("E30").Value = "  -0.05%  "

# Row 31: Filecoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "

# Row 32: Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0547"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "

# Row 34: TrustWalletToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.84%  "

# Row 35: LidoDAOToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.39%  "

# Row 36: ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.689"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "

# Row 37: Aave
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "91.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.49%  "

# Row 38: WEMIXToken
$ws.Range("E38").Value = "  +5.68%  "

# Row 39: Maker
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.329.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.41%  "

# Row 40: VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0192"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.79%  "

# Row 41: ARBITRUM
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.973"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.17%  "

# Row 42: RenderToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.78%  "

# Row 43: HuobiToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.65%  "

# Row 44: InjectiveProtocol
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.35%  "

# Row 45: MXToken
$ws.Range("E45").Value = "  -1.55%  "

# Row 46: Kaspa
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0521"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.37%  "

# Row 47: FraxShare
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.30%  "

# Row 48: RocketPoolETH
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.998.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "

# Row 49: PaxDollar
$ws.Range("E49").Value = "  -0.17%  "

# Row 50: Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0662"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.67%  "

# Row 51: Quant
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "96.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.57%  "

